$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45179 -> 45180, i.e. 2023-09-10 -> 2023-09-11) for every data row (2..200).
for ($row = 2; $row -le 200; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180
}
